$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 787, pushing the old
# rows 787:804 down to 791:808.
$ws.Rows("787:790").Insert()

# New weekly data block (fecha = 2021-09-09, serial 44448) that now
# occupies rows 787:790.
$data = @(
  @(44448, "Primera", 2920, 4500, 5500, 4974, "Chillán", 249),
  @(44448, "Primera", 1100, 5000, 5000, 5000, "Región Metropolitana", 250),
  @(44448, "Segunda", 770, 4000, 4500, 4247, "Chillán", 212),
  @(44448, "Segunda", 380, 4000, 4000, 4000, "Región Metropolitana", 200)
)

$r = 787
foreach ($row in $data) {
    $ws.Range("A$r").Value = 6
    $ws.Range("B$r").Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Range("C$r").Value = "Metropolitana"
    $ws.Range("D$r").Value = $row[0]
    $ws.Range("E$r").Value = 13
    $ws.Range("F$r").Value = 100114013
    $ws.Range("G$r").Value = "Zanahoria"
    $ws.Range("H$r").Value = "Sin especificar"
    $ws.Range("I$r").Value = $row[1]
    $ws.Range("J$r").Value = $row[2]
    $ws.Range("K$r").Value = $row[3]
    $ws.Range("L$r").Value = $row[4]
    $ws.Range("M$r").Value = $row[5]
    $ws.Range("N$r").Value = "$/saco 20 kilos"
    $ws.Range("O$r").Value = $row[6]
    $ws.Range("P$r").Value = $row[7]
    $ws.Range("Q$r").Value = 20
    $ws.Range("R$r").Value = "Hortaliza"
    $r = $r + 1
}
